$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.023.47'
$ws.Range('E2').Value = '  +5.01%  '
$ws.Range('D3').Value = '3.246.68'
$ws.Range('E3').Value = '  +2.65%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''395.48'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').Value = '''108.09'
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('D7').Value = '''0.587'
$ws.Range('E7').Value = '  +6.87%  '
$ws.Range('D8').Value = '3.241.80'
$ws.Range('E8').Value = '  +2.71%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '''0.626'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').Value = '''39.21'
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').Value = '''0.0985'
$ws.Range('E12').Value = '  +12.18%  '
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = '3.759.86'
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').Value = '''8.20'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '3.243.76'
$ws.Range('E17').Value = '  +2.51%  '
$ws.Range('E18').Value = '  -2.71%  '
$ws.Range('D19').Value = '''10.83'
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('D20').Value = '56.820.21'
$ws.Range('E20').Value = '  +5.14%  '
$ws.Range('E21').Value = '  +2.26%  '
$ws.Range('D22').Value = '''0.0000113'
$ws.Range('E22').Value = '  +15.11%  '
$ws.Range('D23').Value = '''13.01'
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('D24').Value = '''295.94'
$ws.Range('E24').Value = '  +8.90%  '
$ws.Range('D25').Value = '''74.31'
$ws.Range('E25').Value = '  +4.51%  '
$ws.Range('D26').Value = '''3.17'
$ws.Range('E26').Value = '  -3.01%  '
$ws.Range('D27').Value = '''27.85'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').Value = '''7.61'
$ws.Range('E28').Value = '  -5.18%  '
$ws.Range('D29').Value = '''7.29'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = '''11.43'
$ws.Range('E32').Value = '  +3.94%  '
$ws.Range('E33').Value = '  -2.67%  '
$ws.Range('D34').Value = '''39.42'
$ws.Range('E34').Value = '  +6.50%  '
$ws.Range('E35').Value = '  -4.79%  '
$ws.Range('E36').Value = '  +1.72%  '
$ws.Range('D37').Value = '''51.62'
$ws.Range('E37').Value = '  +2.21%  '
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('D39').Value = '''0.998'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '''2.91'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('D41').Value = '''135.29'
$ws.Range('E41').Value = '  +3.56%  '
$ws.Range('E42').Value = '  +4.19%  '
$ws.Range('D43').Value = '''17.10'
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '''1.89'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '''3.95'
$ws.Range('E45').Value = '  -3.92%  '
$ws.Range('D46').Value = '''0.282'
$ws.Range('E46').Value = '  -3.49%  '
$ws.Range('D47').Value = '''22.24'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('E48').Value = '  +2.59%  '
$ws.Range('D49').Value = '2.161.52'
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('D50').Value = '''1.99'
$ws.Range('E50').Value = '  +19.04%  '
$ws.Range('E51').Value = '  -3.72%  '
